# The document had four headings that began with a leading "⭐" run
# (rendered in the "Segoe UI Emoji" font, bold) immediately followed by a
# second run holding the heading text itself (with its own leading space),
# e.g.:
#
#   <w:r><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" .../><w:b/><w:bCs/></w:rPr><w:t>⭐</w:t></w:r>
#   <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> 1. Prompt ban đầu</w:t></w:r>
#
# The edit removes the star-emoji run from each of those four headings,
# leaving the following run (and its leading space) untouched, e.g.:
#
#   <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> 1. Prompt ban đầu</w:t></w:r>
#
# A plain Find/Replace-All of the "⭐" glyph with an empty string removes
# just that character (and, because it was the only content of its run,
# Word drops the now-empty "Segoe UI Emoji" run entirely), which reproduces
# the change for all four occurrences in one pass.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "⭐",    # FindText
    $false,  # MatchCase
    $false,  # MatchWholeWord
    $false,  # MatchWildcards
    $false,  # MatchSoundsLike
    $false,  # MatchAllWordForms
    $true,   # Forward
    1,       # Wrap (wdFindContinue)
    $false,  # Format
    "",      # ReplaceWith
    2        # Replace (wdReplaceAll)
)
